# OpenTBS 1.9.1-beta-2014-07-22 : credit features
#
# The canonical-OOXML diff for this revision only renumbers the internal
# chart axis identifiers on the bar3DChart found on slide 3
# (ppt/charts/chart1.xml):
#
#   category axis (c:catAx)  95843456 -> 61990016
#   value    axis (c:valAx)  95844992 -> 61991552
#
# and updates every <c:axId>/<c:crossAx> cross-reference to match. Nothing
# else about the chart (series, data, formatting, 3-D shape, position, ...)
# changes.
#
# Locate the chart shape (slide 3, the "Graphique 3" chart frame) and
# re-stamp its two axis ids through the standard PowerPoint Chart/Axes
# object model, mirroring a user re-touching the chart's axis options.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

$chartShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $candidate = $s.Shapes.Item($i)
    if ($candidate.HasChart) {
        $chartShape = $candidate
        break
    }
}

$chart = $chartShape.Chart

$catAxis = $chart.Axes(1)   # xlCategory -> <c:catAx>
$valAxis = $chart.Axes(2)   # xlValue    -> <c:valAx>

try {
    $catAxis.AxisId = 61990016
    $valAxis.AxisId = 61991552
} catch {
    Write-Host "AxisId reassignment not available: $_"
}
